# Updates the Price (column D) and Volume(1h) (column E) columns of the
# cryptos list per the GitHub Actions refresh commit.
#
# Note: several new Price values (e.g. "249.80", "11.10") are digit
# strings that Excel's input parser would otherwise auto-convert to a
# Number (losing the trailing zero, e.g. "249.80" -> 249.8). Those are
# written with a leading single-quote, the standard Excel "force text"
# prefix, so the cell stays text and the exact original digits are kept.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.457.92"
$ws.Range("D3").Value = "3.356.61"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'249.80"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "'655.45"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  -9.89%  "
$ws.Range("E8").Value = "  -10.47%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -5.91%  "
$ws.Range("D11").Value = "3.349.78"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").Value = "'40.59"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").Value = "97.318.86"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "'6.09"
$ws.Range("E15").Value = "  +7.28%  "
$ws.Range("E16").Value = "  -5.62%  "
$ws.Range("D17").Value = "3.976.41"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'8.67"
$ws.Range("E18").Value = "  +10.64%  "
$ws.Range("D19").Value = "3.355.13"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'0.554"
$ws.Range("E20").Value = "  +25.18%  "
$ws.Range("D21").Value = "'16.89"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "'10.78"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "'504.40"
$ws.Range("E23").Value = "  -4.80%  "
$ws.Range("D24").Value = "'3.35"
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("E25").Value = "  -6.46%  "
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "'96.79"
$ws.Range("E27").Value = "  -5.93%  "
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").Value = "3.537.04"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "'0.995"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "'11.10"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +19.42%  "
$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.553"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("D37").Value = "'28.66"
$ws.Range("D38").Value = "'7.75"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "'1.46"
$ws.Range("E39").Value = "  +9.74%  "
$ws.Range("D40").Value = "'519.35"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").Value = "'24.63"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "'8.84"
$ws.Range("E45").Value = "  +12.61%  "
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "'3.68"
$ws.Range("E47").Value = "  -6.37%  "
$ws.Range("D48").Value = "'5.59"
$ws.Range("E48").Value = "  +8.11%  "
$ws.Range("E49").Value = "  +7.21%  "
$ws.Range("D50").Value = "'53.37"
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("E51").Value = "  -6.29%  "
